$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the active selection to G12
$ws.Range("G12").Select()

# Add new column S values for rows 10-20 (match neighboring cells: default/no explicit style)
$ws.Range("S10").Value = 0
$ws.Range("S10").Style = "Normal"
$ws.Range("S11").Value = 1
$ws.Range("S11").Style = "Normal"
$ws.Range("S12").Value = 2.391048198037748
$ws.Range("S12").Style = "Normal"
$ws.Range("S13").Value = 1.529901383525567
$ws.Range("S13").Style = "Normal"
$ws.Range("S14").Value = -0.06158374816301929
$ws.Range("S14").Style = "Normal"
$ws.Range("S15").Value = -0.06158374816301929
$ws.Range("S15").Style = "Normal"
$ws.Range("S16").Value = 1.160087863853805
$ws.Range("S16").Style = "Normal"
$ws.Range("S17").Value = 1.231283191936991
$ws.Range("S17").Style = "Normal"
$ws.Range("S18").Value = 1.637451109109422
$ws.Range("S18").Style = "Normal"
$ws.Range("S19").Value = 3.146844050165727
$ws.Range("S19").Style = "Normal"
$ws.Range("S20").Value = 5.680878344533673
$ws.Range("S20").Style = "Normal"

# Update row 31 niter
$ws.Range("S31").Value = 558
$ws.Range("S31").Style = "Normal"

# Update row 32 fit time - update B32:R32 and add S32
$ws.Range("B32").Value = 101.809488558
$ws.Range("C32").Value = 161.817944822
$ws.Range("D32").Value = 116.464285007
$ws.Range("E32").Value = 169.800406674
$ws.Range("F32").Value = 134.699069528
$ws.Range("G32").Value = 473.537694772
$ws.Range("H32").Value = 67.42791133499986
$ws.Range("I32").Value = 58.54798157100004
$ws.Range("J32").Value = 100.419380647
$ws.Range("K32").Value = 177.740926278
$ws.Range("L32").Value = 136.3318007939999
$ws.Range("M32").Value = 528.9777407750003
$ws.Range("N32").Value = 90.44575723599974
$ws.Range("O32").Value = 166.885690393
$ws.Range("P32").Value = 100.5118735900001
$ws.Range("Q32").Value = 201.939431191
$ws.Range("R32").Value = 272.4740447959998
$ws.Range("S32").Value = 220.9807230359997
$ws.Range("S32").Style = "Normal"

# Row 33 chi-square
$ws.Range("S33").Value = 6.679512131347471
$ws.Range("S33").Style = "Normal"

# Row 34 ndf
$ws.Range("S34").Value = 7
$ws.Range("S34").Style = "Normal"

# Row 35 chi2-per-ndf
$ws.Range("S35").Value = 0.9542160187639245
$ws.Range("S35").Style = "Normal"

# Row 37 number lineups
$ws.Range("S37").Value = 640
$ws.Range("S37").Style = "Normal"

# Row 38 number TA lineups
$ws.Range("S38").Value = 311
$ws.Range("S38").Style = "Normal"

# Row 39 number TP lineups
$ws.Range("S39").Value = 329
$ws.Range("S39").Style = "Normal"

# Row 40 correctID
$ws.Range("S40").Value = 0.547112462006079
$ws.Range("S40").Style = "Normal"

# Row 41 falseID
$ws.Range("S41").Value = 0.05627009646302251
$ws.Range("S41").Style = "Normal"

# Row 42 d-prime
$ws.Range("S42").Value = 1.705247644647274
$ws.Range("S42").Style = "Normal"

# Row 43 pAUC
$ws.Range("S43").Value = 0.021403240957212
$ws.Range("S43").Style = "Normal"
